# Adds the title-slide text ("Elevator UML" / "By the us") that turns the
# blank title/subtitle placeholders on slide 1 into a proper title page, and
# shrinks the now-autofitted subtitle placeholder to its new rendered height.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 2 = "Title 1" (ctrTitle placeholder), Shape 3 = "Subtitle 2" (subTitle placeholder)
$title    = $s.Shapes.Item(2)
$subtitle = $s.Shapes.Item(3)

# --- Title: "Elevator UML" ---
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Elevator UML"
# Keep the run's language the same as the rest of the deck (en-AU) instead of
# the engine's default en-US.
$titleRange.LanguageID = "en-AU"

# --- Subtitle: "By the us" (14pt) ---
$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.Text = "By the us"
$subtitleRange.Font.Size = 14
$subtitleRange.LanguageID = "en-AU"

# The subtitle placeholder auto-fits its box to the (now non-empty) text,
# shrinking the height from 1655762 EMU (130.375pt) to 1251889 EMU.
$subtitle.Height = 98.57397700787402
